$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (used in A1)
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 10:55"

# Row 19 (Banglades) - updated stats
$ws.Range("B19").Value = 234889
$ws.Range("C19").Value = 2695
$ws.Range("D19").Value = 132960
$ws.Range("E19").Value = 98846
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = 3083

# Row 27 (Indonesia) - updated stats
$ws.Range("B27").Value = 106336
$ws.Range("C27").Value = 1904
$ws.Range("D27").Value = 64292
$ws.Range("E27").Value = 36986
$ws.Range("G27").Value = 83
$ws.Range("H27").Value = 5058

# Row 45 (Singapur) - updated stats
$ws.Range("B45").Value = 51809
$ws.Range("C45").Value = 278
$ws.Range("E45").Value = 5684

# Row 49 (Polonia) - updated stats
$ws.Range("B49").Value = 45031
$ws.Range("C49").Value = 615
$ws.Range("D49").Value = 33643
$ws.Range("E49").Value = 9679
$ws.Range("G49").Value = 15
$ws.Range("H49").Value = 1709

# Row 66 (Austria) - updated stats
$ws.Range("B66").Value = 20955
$ws.Range("C66").Value = 105
$ws.Range("D66").Value = 18628
$ws.Range("E66").Value = 1609
$ws.Range("G66").Value = 2
$ws.Range("H66").Value = 718

# Rows 73/74 (Chequia / El Salvador) swap order and El Salvador gets new stats.
# El Salvador now ranks above Chequia, so it takes row 73 with its updated
# numbers, while Chequia (unchanged numbers) moves down to row 74.
$ws.Range("A73").Value = "El Salvador"
$ws.Range("B73").Value = 16230
$ws.Range("C73").Value = 389
$ws.Range("D73").Value = 8206
$ws.Range("E73").Value = 7585
$ws.Range("G73").Value = 9
$ws.Range("H73").Value = 439

$ws.Range("A74").Value = "Chequia"
$ws.Range("B74").Value = 16093
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 11429
$ws.Range("E74").Value = 4290
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 374

# Row 124 (Eslovaquia) - updated stats
$ws.Range("B124").Value = 2265
$ws.Range("C124").Value = 20
$ws.Range("D124").Value = 1675
$ws.Range("E124").Value = 562
